# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header row (bold, centered, thin box border)
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Fill every data row (2 through 45) with the team's season record
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 65
    $ws.Cells.Item($r, 31).Value = 97
    $ws.Cells.Item($r, 32).Value = 0
}
